# Insert two new data rows (weekly Fruta/Hortalizas update) above the existing
# row 27, shifting all subsequent rows down by two (old row N -> new row N+2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 27; formatting (date style on column D) is
# inherited from the row above, matching the rest of the sheet.
$ws.Range("A27:A28").EntireRow.Insert()

# --- New row 27: Larry Ann, Primera, Region de O'Higgins ---
$ws.Cells.Item(27, 1).Value = 9
$ws.Cells.Item(27, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value = "Metropolitana"
$ws.Cells.Item(27, 4).Value = 44607
$ws.Cells.Item(27, 5).Value = 13
$ws.Cells.Item(27, 6).Value = "Fruta"
$ws.Cells.Item(27, 7).Value = 100103
$ws.Cells.Item(27, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(27, 9).Value = 100103002
$ws.Cells.Item(27, 10).Value = "Ciruela"
$ws.Cells.Item(27, 11).Value = "Larry Ann"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 120
$ws.Cells.Item(27, 14).Value = 12500
$ws.Cells.Item(27, 15).Value = 12500
$ws.Cells.Item(27, 16).Value = 12500
$ws.Cells.Item(27, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(27, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(27, 19).Value = 694
$ws.Cells.Item(27, 20).Value = 18

# --- New row 28: Larry Ann, Segunda, Region de O'Higgins ---
$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44607
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100103
$ws.Cells.Item(28, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(28, 9).Value = 100103002
$ws.Cells.Item(28, 10).Value = "Ciruela"
$ws.Cells.Item(28, 11).Value = "Larry Ann"
$ws.Cells.Item(28, 12).Value = "Segunda"
$ws.Cells.Item(28, 13).Value = 150
$ws.Cells.Item(28, 14).Value = 10000
$ws.Cells.Item(28, 15).Value = 10000
$ws.Cells.Item(28, 16).Value = 10000
$ws.Cells.Item(28, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(28, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(28, 19).Value = 556
$ws.Cells.Item(28, 20).Value = 18
